$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$tooltipUrl = "https://fiori.jioconnect.com/sap/bc/ui5_ui5/sap/zehys_dashboard/javascript:void(0);"

# --- New data rows 70-73 (order batches received 11-13 Apr 2021) ---------
# Copy the row-69 number formats down first so A/B/C/D/E pick up the usual
# "S.No / Order ID / Topup / Order / Creation Date" look before we fill them.
$ws.Range("A69:E69").Copy()
$ws.Range("A70:E73").PasteSpecial(-4122)

$ws.Range("A70").Value = 67
$ws.Range("B70").Value = 26718802
$ws.Range("C70").Value = 18992
$ws.Range("D70").Value = 17999.23
$ws.Range("E70").Value = 44297

$ws.Range("A71").Value = 68
$ws.Range("B71").Value = 26715231
$ws.Range("C71").Value = 31654
$ws.Range("D71").Value = 30000.05
$ws.Range("E71").Value = 44297

$ws.Range("A72").Value = 69
$ws.Range("B72").Value = 26724015
$ws.Range("C72").Value = 84410
$ws.Range("D72").Value = 79999.13
$ws.Range("E72").Value = 44298

$ws.Range("A73").Value = 70
$ws.Range("B73").Value = 26744506
$ws.Range("C73").Value = 105513
$ws.Range("D73").Value = 100000.16
$ws.Range("E73").Value = 44299

# F:I already carry the shared IF(B="","",...) formula definition inherited
# from the block above, but the cached results need to be forced to
# recompute now that B:D actually hold values.
foreach ($r in 70..73) {
    $ws.Range("F$r").Formula = "=IF(B$r=`"`",`"`",C$r-D$r)"
    $ws.Range("G$r").Formula = "=IF(B$r=`"`",`"`",F$r/D$r*100)"
    $ws.Range("H$r").Formula = "=IF(B$r=`"`",`"`",D$r*1.04)"
    $ws.Range("I$r").Formula = "=IF(B$r=`"`",`"`",C$r-H$r)"
}

# Order-ID hyperlinks, matching the ones already on B3:B9. Adding a
# hyperlink re-styles the cell with the blue/underlined "Hyperlink" look,
# but B3:B9 show this workbook does NOT use that look for linked order
# IDs, so the plain number formatting is restored right after.
$ws.Hyperlinks.Add($ws.Range("B70"), $tooltipUrl, "", $tooltipUrl, "26718802")
$ws.Hyperlinks.Add($ws.Range("B71"), $tooltipUrl, "", $tooltipUrl, "26715231")
$ws.Hyperlinks.Add($ws.Range("B72"), $tooltipUrl, "", $tooltipUrl, "26724015")
$ws.Hyperlinks.Add($ws.Range("B73"), $tooltipUrl, "", $tooltipUrl, "26744506")

$ws.Range("B69").Copy()
$ws.Range("B70:B73").PasteSpecial(-4122)

# --- Row 74: keep it blank, but give it the usual S.No column formatting -
$ws.Range("A69").Copy()
$ws.Range("A74").PasteSpecial(-4122)

# --- Rows 83-90: extend the blank Topup/Order/Date formatting down -------
$ws.Range("C82:D82").Copy()
$ws.Range("C83:D85").PasteSpecial(-4122)

$ws.Range("C82:E82").Copy()
$ws.Range("C86:E90").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- View state: scrolled a bit further down, selection moved to J65 -----
$excel.Goto($ws.Range("A57"), $false)
$ws.Range("J65").Select()

$excel.Calculate()
